$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Rating") gets the text value "9" written into every row
# (header row 1 through the last data row, 101). Copy a cell that already
# holds the text "9" (A10, the "Sl. No" for the 9th book) and paste its
# value into each D cell so the result stays a text/shared-string "9"
# rather than becoming a numeric 9, and no new cell style is introduced.
$ws.Range("A10").Copy()
for ($r = 1; $r -le 101; $r++) {
    $ws.Cells.Item($r, 4).PasteSpecial(-4163)
}
